$d = $word.ActiveDocument

$d.Content.Find.Execute("32+63=", $true, $false, $false, $false, $false, $true, 1, $false, "28-0=", 2) | Out-Null
$d.Content.Find.Execute("83-58=", $true, $false, $false, $false, $false, $true, 1, $false, "80-56=", 2) | Out-Null
$d.Content.Find.Execute("39-15=", $true, $false, $false, $false, $false, $true, 1, $false, "49+33=", 2) | Out-Null
$d.Content.Find.Execute("74-44=", $true, $false, $false, $false, $false, $true, 1, $false, "9+2=", 2) | Out-Null
$d.Content.Find.Execute("27-12=", $true, $false, $false, $false, $false, $true, 1, $false, "4+28=", 2) | Out-Null
$d.Content.Find.Execute("78-58=", $true, $false, $false, $false, $false, $true, 1, $false, "55+8=", 2) | Out-Null
$d.Content.Find.Execute("19-15=", $true, $false, $false, $false, $false, $true, 1, $false, "95-63=", 2) | Out-Null
$d.Content.Find.Execute("70-18=", $true, $false, $false, $false, $false, $true, 1, $false, "46+48=", 2) | Out-Null
$d.Content.Find.Execute("49-28=", $true, $false, $false, $false, $false, $true, 1, $false, "60-59=", 2) | Out-Null
$d.Content.Find.Execute("83-61=", $true, $false, $false, $false, $false, $true, 1, $false, "34+60=", 2) | Out-Null
$d.Content.Find.Execute("49-7=", $true, $false, $false, $false, $false, $true, 1, $false, "96-1=", 2) | Out-Null
$d.Content.Find.Execute("56-16=", $true, $false, $false, $false, $false, $true, 1, $false, "59-32=", 2) | Out-Null
$d.Content.Find.Execute("81+13=", $true, $false, $false, $false, $false, $true, 1, $false, "85-71=", 2) | Out-Null
$d.Content.Find.Execute("0+34=", $true, $false, $false, $false, $false, $true, 1, $false, "62+35=", 2) | Out-Null
$d.Content.Find.Execute("21+69=", $true, $false, $false, $false, $false, $true, 1, $false, "18+68=", 2) | Out-Null
$d.Content.Find.Execute("14+73=", $true, $false, $false, $false, $false, $true, 1, $false, "70-35=", 2) | Out-Null
$d.Content.Find.Execute("15+63=", $true, $false, $false, $false, $false, $true, 1, $false, "87-52=", 2) | Out-Null
$d.Content.Find.Execute("15+74=", $true, $false, $false, $false, $false, $true, 1, $false, "28+20=", 2) | Out-Null
$d.Content.Find.Execute("90+6=", $true, $false, $false, $false, $false, $true, 1, $false, "74-22=", 2) | Out-Null
$d.Content.Find.Execute("59+14=", $true, $false, $false, $false, $false, $true, 1, $false, "0+83=", 2) | Out-Null
$d.Content.Find.Execute("20+20=", $true, $false, $false, $false, $false, $true, 1, $false, "66-10=", 2) | Out-Null
$d.Content.Find.Execute("45+47=", $true, $false, $false, $false, $false, $true, 1, $false, "3+41=", 2) | Out-Null
$d.Content.Find.Execute("77-11=", $true, $false, $false, $false, $false, $true, 1, $false, "65+2=", 2) | Out-Null
$d.Content.Find.Execute("45-35=", $true, $false, $false, $false, $false, $true, 1, $false, "98-66=", 2) | Out-Null
$d.Content.Find.Execute("79-76=", $true, $false, $false, $false, $false, $true, 1, $false, "83-72=", 2) | Out-Null
$d.Content.Find.Execute("76+9=", $true, $false, $false, $false, $false, $true, 1, $false, "3+0=", 2) | Out-Null
$d.Content.Find.Execute("14+78=", $true, $false, $false, $false, $false, $true, 1, $false, "97-12=", 2) | Out-Null
$d.Content.Find.Execute("8+17=", $true, $false, $false, $false, $false, $true, 1, $false, "3-3=", 2) | Out-Null
$d.Content.Find.Execute("9+34=", $true, $false, $false, $false, $false, $true, 1, $false, "6+2=", 2) | Out-Null
$d.Content.Find.Execute("37-5=", $true, $false, $false, $false, $false, $true, 1, $false, "11-5=", 2) | Out-Null
$d.Content.Find.Execute("67-61=", $true, $false, $false, $false, $false, $true, 1, $false, "39+9=", 2) | Out-Null
$d.Content.Find.Execute("0+97=", $true, $false, $false, $false, $false, $true, 1, $false, "73-28=", 2) | Out-Null
$d.Content.Find.Execute("56-6=", $true, $false, $false, $false, $false, $true, 1, $false, "75-4=", 2) | Out-Null
$d.Content.Find.Execute("67-16=", $true, $false, $false, $false, $false, $true, 1, $false, "6+57=", 2) | Out-Null
$d.Content.Find.Execute("48-12=", $true, $false, $false, $false, $false, $true, 1, $false, "46+3=", 2) | Out-Null
$d.Content.Find.Execute("1+98=", $true, $false, $false, $false, $false, $true, 1, $false, "11+83=", 2) | Out-Null
$d.Content.Find.Execute("37-22=", $true, $false, $false, $false, $false, $true, 1, $false, "33-31=", 2) | Out-Null
$d.Content.Find.Execute("83-20=", $true, $false, $false, $false, $false, $true, 1, $false, "89-23=", 2) | Out-Null
$d.Content.Find.Execute("40-35=", $true, $false, $false, $false, $false, $true, 1, $false, "36+41=", 2) | Out-Null
$d.Content.Find.Execute("31-20=", $true, $false, $false, $false, $false, $true, 1, $false, "12+50=", 2) | Out-Null
$d.Content.Find.Execute("58-52=", $true, $false, $false, $false, $false, $true, 1, $false, "6+14=", 2) | Out-Null
$d.Content.Find.Execute("48+48=", $true, $false, $false, $false, $false, $true, 1, $false, "92-79=", 2) | Out-Null
$d.Content.Find.Execute("14+15=", $true, $false, $false, $false, $false, $true, 1, $false, "66-19=", 2) | Out-Null
$d.Content.Find.Execute("49-49=", $true, $false, $false, $false, $false, $true, 1, $false, "97-37=", 2) | Out-Null
$d.Content.Find.Execute("37+62=", $true, $false, $false, $false, $false, $true, 1, $false, "27+20=", 2) | Out-Null
$d.Content.Find.Execute("51+48=", $true, $false, $false, $false, $false, $true, 1, $false, "42+31=", 2) | Out-Null
$d.Content.Find.Execute("37-13=", $true, $false, $false, $false, $false, $true, 1, $false, "77-24=", 2) | Out-Null
$d.Content.Find.Execute("73-69=", $true, $false, $false, $false, $false, $true, 1, $false, "21+23=", 2) | Out-Null
$d.Content.Find.Execute("74-55=", $true, $false, $false, $false, $false, $true, 1, $false, "75+13=", 2) | Out-Null
$d.Content.Find.Execute("21-9=", $true, $false, $false, $false, $false, $true, 1, $false, "94-43=", 2) | Out-Null
$d.Content.Find.Execute("63-50=", $true, $false, $false, $false, $false, $true, 1, $false, "46+44=", 2) | Out-Null
$d.Content.Find.Execute("42-18=", $true, $false, $false, $false, $false, $true, 1, $false, "57+42=", 2) | Out-Null
$d.Content.Find.Execute("96-2=", $true, $false, $false, $false, $false, $true, 1, $false, "58-58=", 2) | Out-Null
$d.Content.Find.Execute("56-39=", $true, $false, $false, $false, $false, $true, 1, $false, "84-0=", 2) | Out-Null
$d.Content.Find.Execute("84-14=", $true, $false, $false, $false, $false, $true, 1, $false, "73-0=", 2) | Out-Null
$d.Content.Find.Execute("26+73=", $true, $false, $false, $false, $false, $true, 1, $false, "70+29=", 2) | Out-Null
$d.Content.Find.Execute("65-45=", $true, $false, $false, $false, $false, $true, 1, $false, "48+31=", 2) | Out-Null
$d.Content.Find.Execute("25+24=", $true, $false, $false, $false, $false, $true, 1, $false, "52-46=", 2) | Out-Null
$d.Content.Find.Execute("53+39=", $true, $false, $false, $false, $false, $true, 1, $false, "90-31=", 2) | Out-Null
$d.Content.Find.Execute("75-34=", $true, $false, $false, $false, $false, $true, 1, $false, "70-39=", 2) | Out-Null
$d.Content.Find.Execute("96-36=", $true, $false, $false, $false, $false, $true, 1, $false, "75-52=", 2) | Out-Null
$d.Content.Find.Execute("48+6=", $true, $false, $false, $false, $false, $true, 1, $false, "97-43=", 2) | Out-Null
$d.Content.Find.Execute("99-64=", $true, $false, $false, $false, $false, $true, 1, $false, "25+0=", 2) | Out-Null
$d.Content.Find.Execute("32+31=", $true, $false, $false, $false, $false, $true, 1, $false, "21+16=", 2) | Out-Null
$d.Content.Find.Execute("99-78=", $true, $false, $false, $false, $false, $true, 1, $false, "90-68=", 2) | Out-Null
$d.Content.Find.Execute("52-18=", $true, $false, $false, $false, $false, $true, 1, $false, "43-11=", 2) | Out-Null
$d.Content.Find.Execute("35-33=", $true, $false, $false, $false, $false, $true, 1, $false, "98-23=", 2) | Out-Null
$d.Content.Find.Execute("80-61=", $true, $false, $false, $false, $false, $true, 1, $false, "28+28=", 2) | Out-Null
$d.Content.Find.Execute("88-43=", $true, $false, $false, $false, $false, $true, 1, $false, "22+11=", 2) | Out-Null
$d.Content.Find.Execute("31+68=", $true, $false, $false, $false, $false, $true, 1, $false, "45-1=", 2) | Out-Null
$d.Content.Find.Execute("72-26=", $true, $false, $false, $false, $false, $true, 1, $false, "1+62=", 2) | Out-Null
$d.Content.Find.Execute("87-31=", $true, $false, $false, $false, $false, $true, 1, $false, "56-38=", 2) | Out-Null
$d.Content.Find.Execute("43+3=", $true, $false, $false, $false, $false, $true, 1, $false, "92-90=", 2) | Out-Null
$d.Content.Find.Execute("35+31=", $true, $false, $false, $false, $false, $true, 1, $false, "92-40=", 2) | Out-Null
$d.Content.Find.Execute("76+7=", $true, $false, $false, $false, $false, $true, 1, $false, "11+68=", 2) | Out-Null
$d.Content.Find.Execute("50+32=", $true, $false, $false, $false, $false, $true, 1, $false, "76-3=", 2) | Out-Null
$d.Content.Find.Execute("23+7=", $true, $false, $false, $false, $false, $true, 1, $false, "43-27=", 2) | Out-Null
$d.Content.Find.Execute("74+6=", $true, $false, $false, $false, $false, $true, 1, $false, "27-26=", 2) | Out-Null
$d.Content.Find.Execute("59-45=", $true, $false, $false, $false, $false, $true, 1, $false, "11+42=", 2) | Out-Null
$d.Content.Find.Execute("17+17=", $true, $false, $false, $false, $false, $true, 1, $false, "43+7=", 2) | Out-Null
$d.Content.Find.Execute("89-52=", $true, $false, $false, $false, $false, $true, 1, $false, "69-20=", 2) | Out-Null
$d.Content.Find.Execute("5+70=", $true, $false, $false, $false, $false, $true, 1, $false, "65+24=", 2) | Out-Null
$d.Content.Find.Execute("7+30=", $true, $false, $false, $false, $false, $true, 1, $false, "64+25=", 2) | Out-Null
$d.Content.Find.Execute("27+47=", $true, $false, $false, $false, $false, $true, 1, $false, "21-19=", 2) | Out-Null
$d.Content.Find.Execute("62+5=", $true, $false, $false, $false, $false, $true, 1, $false, "71+0=", 2) | Out-Null
$d.Content.Find.Execute("65-46=", $true, $false, $false, $false, $false, $true, 1, $false, "80-8=", 2) | Out-Null
$d.Content.Find.Execute("1+85=", $true, $false, $false, $false, $false, $true, 1, $false, "51+26=", 2) | Out-Null
$d.Content.Find.Execute("86-6=", $true, $false, $false, $false, $false, $true, 1, $false, "59-59=", 2) | Out-Null
$d.Content.Find.Execute("92-67=", $true, $false, $false, $false, $false, $true, 1, $false, "71-19=", 2) | Out-Null
$d.Content.Find.Execute("13+32=", $true, $false, $false, $false, $false, $true, 1, $false, "52+33=", 2) | Out-Null
$d.Content.Find.Execute("43-35=", $true, $false, $false, $false, $false, $true, 1, $false, "52+19=", 2) | Out-Null
$d.Content.Find.Execute("80-3=", $true, $false, $false, $false, $false, $true, 1, $false, "31+30=", 2) | Out-Null
$d.Content.Find.Execute("68-55=", $true, $false, $false, $false, $false, $true, 1, $false, "37+23=", 2) | Out-Null
$d.Content.Find.Execute("82+2=", $true, $false, $false, $false, $false, $true, 1, $false, "65-12=", 2) | Out-Null
$d.Content.Find.Execute("19+72=", $true, $false, $false, $false, $false, $true, 1, $false, "80-45=", 2) | Out-Null
$d.Content.Find.Execute("78-15=", $true, $false, $false, $false, $false, $true, 1, $false, "72-71=", 2) | Out-Null
$d.Content.Find.Execute("13+37=", $true, $false, $false, $false, $false, $true, 1, $false, "38+51=", 2) | Out-Null
$d.Content.Find.Execute("89+10=", $true, $false, $false, $false, $false, $true, 1, $false, "96-25=", 2) | Out-Null
$d.Content.Find.Execute("38+10=", $true, $false, $false, $false, $false, $true, 1, $false, "93-56=", 2) | Out-Null
$d.Content.Find.Execute("67-13=", $true, $false, $false, $false, $false, $true, 1, $false, "61+18=", 2) | Out-Null
